$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.258.87"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "3.058.75"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "548.94"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").Value = "139.80"
$ws.Range("E6").Value = "  +4.31%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.052.90"
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "6.40"
$ws.Range("E10").Value = "  +5.46%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "3.550.73"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").Value = "63.299.88"
$ws.Range("E16").Value = "  +3.36%  "
$ws.Range("D17").Value = "3.056.21"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").Value = "6.73"
$ws.Range("E19").Value = "  +1.83%  "
$ws.Range("D20").Value = "481.17"
$ws.Range("E20").Value = "  +3.39%  "
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "7.21"
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("D24").Value = "80.63"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "12.51"
$ws.Range("E25").Value = "  +3.64%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").Value = "7.92"
$ws.Range("E29").Value = "  +4.44%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("D31").Value = "25.97"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  +5.95%  "
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("D35").Value = "55.35"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "465.06"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").Value = "0.0817"
$ws.Range("E38").Value = "  +3.81%  "
$ws.Range("D39").Value = "0.0396"
$ws.Range("E39").Value = "  +2.74%  "
$ws.Range("D40").Value = "3.064.79"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "28.31"
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D47").Value = "2.05"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").Value = "116.93"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("E51").Value = "  +1.95%  "
